# refactor: the logic class structure in xml files
#
# The sheet gained a new "Force" option row, inserted right after the
# existing "Cache" row (old row 7) and before the former "Upload" row
# (old row 8), which shifts down to become the new row 9. Everything
# below (the "Ref" row and the sample data row) shifts down by one row
# as a consequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8 - this pushes the old row 8
# ("Upload") and everything after it down by one, and Excel automatically
# extends the data validation sqref ranges (A6:A8 -> A6:A9,
# B7:J8 -> B7:J9) and the sheet dimension to match.
$ws.Rows.Item(8).Insert()

# Copy the formatting of the row above (row 7, "Cache") onto the new
# row 8 so it picks up the same per-cell styles (label cell + boolean
# cells) used by the other option rows.
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)

# Fill in the new row's content: label "Force" in column A and FALSE
# booleans across B:I, matching the pattern of the other rows.
$ws.Range("A8").Value = "Force"
$ws.Range("B8:I8").Value = $false

# Re-freeze the panes so the split/frozen row follows the inserted row
# (the header block is now 10 rows tall instead of 9) and restore the
# selection to the newly added row.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A11").Select()
$win.FreezePanes = $true
$ws.Range("A9").Select()
